# Insert 3 new weekly price rows for "Vega Monumental Concepción - Plátano"
# at row 677, pushing the existing rows (677-726) down to (680-729).
# Excel's Rows.Insert() carries formatting (e.g. the date style on column D)
# down onto the newly inserted blank rows automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(677).Resize(3).Insert()

# --- Row 677: Maduro ---
$ws.Range("A677").Value = 11
$ws.Range("B677").Value = "Vega Monumental Concepción"
$ws.Range("C677").Value = "Bíobío"
$ws.Range("D677").Value = 45013
$ws.Range("E677").Value = 8
$ws.Range("F677").Value = "Fruta"
$ws.Range("G677").Value = 100108
$ws.Range("H677").Value = "Tropicales y subtropicales"
$ws.Range("I677").Value = 100108006
$ws.Range("J677").Value = "Plátano"
$ws.Range("K677").Value = "Sin especificar"
$ws.Range("L677").Value = "Maduro"
$ws.Range("M677").Value = 100
$ws.Range("N677").Value = 20000
$ws.Range("O677").Value = 20000
$ws.Range("P677").Value = 20000
$ws.Range("Q677").Value = "$/caja 20 kilos"
$ws.Range("R677").Value = "Ecuador"
$ws.Range("S677").Value = 1000
$ws.Range("T677").Value = 20

# --- Row 678: Pintón ---
$ws.Range("A678").Value = 11
$ws.Range("B678").Value = "Vega Monumental Concepción"
$ws.Range("C678").Value = "Bíobío"
$ws.Range("D678").Value = 45013
$ws.Range("E678").Value = 8
$ws.Range("F678").Value = "Fruta"
$ws.Range("G678").Value = 100108
$ws.Range("H678").Value = "Tropicales y subtropicales"
$ws.Range("I678").Value = 100108006
$ws.Range("J678").Value = "Plátano"
$ws.Range("K678").Value = "Sin especificar"
$ws.Range("L678").Value = "Pintón"
$ws.Range("M678").Value = 400
$ws.Range("N678").Value = 21000
$ws.Range("O678").Value = 21000
$ws.Range("P678").Value = 21000
$ws.Range("Q678").Value = "$/caja 20 kilos"
$ws.Range("R678").Value = "Ecuador"
$ws.Range("S678").Value = 1050
$ws.Range("T678").Value = 20

# --- Row 679: Primera Pintón ---
$ws.Range("A679").Value = 11
$ws.Range("B679").Value = "Vega Monumental Concepción"
$ws.Range("C679").Value = "Bíobío"
$ws.Range("D679").Value = 45013
$ws.Range("E679").Value = 8
$ws.Range("F679").Value = "Fruta"
$ws.Range("G679").Value = 100108
$ws.Range("H679").Value = "Tropicales y subtropicales"
$ws.Range("I679").Value = 100108006
$ws.Range("J679").Value = "Plátano"
$ws.Range("K679").Value = "Sin especificar"
$ws.Range("L679").Value = "Primera Pintón"
$ws.Range("M679").Value = 400
$ws.Range("N679").Value = 24000
$ws.Range("O679").Value = 24000
$ws.Range("P679").Value = 24000
$ws.Range("Q679").Value = "$/caja 20 kilos"
$ws.Range("R679").Value = "Ecuador"
$ws.Range("S679").Value = 1200
$ws.Range("T679").Value = 20
